$wb = $excel.ActiveWorkbook

# --- Monsters sheet (sheet3): stat tweaks ---
$wsMonsters = $wb.Worksheets.Item("Monsters")

$wsMonsters.Range("I2").Value = 5

$wsMonsters.Range("H3").Value = 80
$wsMonsters.Range("I3").Value = 6

$wsMonsters.Range("H4").Value = 100
$wsMonsters.Range("I4").Value = 7

$wsMonsters.Range("G5").Value = 350
$wsMonsters.Range("H5").Value = 350
$wsMonsters.Range("I5").Value = 60

$wsMonsters.Range("G6").Value = 350
$wsMonsters.Range("H6").Value = 350
$wsMonsters.Range("I6").Value = 60

$wsMonsters.Range("B7").Value = 80
$wsMonsters.Range("G7").Value = 600
$wsMonsters.Range("H7").Value = 600
$wsMonsters.Range("I7").Value = 150

# --- Sheet2 (sheet6): level table tweak ---
$wsLevel = $wb.Worksheets.Item("Sheet2")
$wsLevel.Range("G3").Value = 400

# --- Selections / active sheet / active cell updates ---
$wsCharacters = $wb.Worksheets.Item("Characters")
$wsCharacters.Range("H2").Select()

$wsLevel.Range("B12").Select()

$wsMonsters.Activate()
$wsMonsters.Range("I5").Select()
